$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.609.22"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.584.68"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.33"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.29"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.584.27"
$ws.Range("E7").Value = "  +1.47%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("E10").Value = "  +0.64%  "
$ws.Range("E11").Value = "  +5.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.392"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.196.13"
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("E14").Value = "  +3.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000186"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.588.70"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.743.12"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.00"
$ws.Range("E19").Value = "  -3.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.59"
$ws.Range("E20").Value = "  +2.15%  "
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "395.54"
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("E23").Value = "  +3.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.733.35"
$ws.Range("E24").Value = "  +1.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.26"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +2.96%  "
$ws.Range("E28").Value = "  +5.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.66"
$ws.Range("E29").Value = "  +31.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.39"
$ws.Range("E30").Value = "  +4.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.59"
$ws.Range("E31").Value = "  +5.13%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.590.32"
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.46"
$ws.Range("E34").Value = "  +2.99%  "
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.42"
$ws.Range("E37").Value = "  +9.03%  "
$ws.Range("E38").Value = "  +4.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.06"
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "169.39"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  +4.70%  "
$ws.Range("E42").Value = "  +2.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.99"
$ws.Range("E43").Value = "  +3.73%  "
$ws.Range("E44").Value = "  +8.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.12"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("E46").Value = "  +2.91%  "
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.03"
$ws.Range("E49").Value = "  +3.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.451.99"
$ws.Range("E50").Value = "  +2.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "317.09"
$ws.Range("E51").Value = "  +5.21%  "
